$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting down from row 600 to 614 (copies cell styles A:H)
$ws.Range("A600:H600").Copy() | Out-Null
$ws.Range("A601:H614").PasteSpecial(-4122) | Out-Null

# Blank "Localisation douleur" cells need the blank-cell style (copy from G598)
$ws.Range("G598").Copy() | Out-Null
$ws.Range("G604").PasteSpecial(-4122) | Out-Null
$ws.Range("G598").Copy() | Out-Null
$ws.Range("G608").PasteSpecial(-4122) | Out-Null
$ws.Range("G598").Copy() | Out-Null
$ws.Range("G609").PasteSpecial(-4122) | Out-Null
$ws.Range("G598").Copy() | Out-Null
$ws.Range("G612").PasteSpecial(-4122) | Out-Null
$ws.Range("G598").Copy() | Out-Null
$ws.Range("G613").PasteSpecial(-4122) | Out-Null

$ws.Range("A601").Value = 45986
$ws.Range("B601").Value = 'Yoan Zouma'
$ws.Range("C601").Value = 70
$ws.Range("D601").Value = 3
$ws.Range("E601").Value = 6
$ws.Range("F601").Value = 5
$ws.Range("G601").Value = 'Cheville,ischio'
$ws.Range("H601").Value = 10
$ws.Range("I601").Formula = "=C601*D601"

$ws.Range("A602").Value = 45986
$ws.Range("B602").Value = 'Yoann Martelat'
$ws.Range("C602").Value = 70
$ws.Range("D602").Value = 4
$ws.Range("E602").Value = 4
$ws.Range("F602").Value = 6
$ws.Range("G602").Value = 'Genou'
$ws.Range("H602").Value = 5
$ws.Range("I602").Formula = "=C602*D602"

$ws.Range("A603").Value = 45986
$ws.Range("B603").Value = 'Kamal Bafounta'
$ws.Range("C603").Value = 70
$ws.Range("D603").Value = 4
$ws.Range("E603").Value = 1
$ws.Range("F603").Value = 1
$ws.Range("G603").Value = 'Genou'
$ws.Range("H603").Value = 8
$ws.Range("I603").Formula = "=C603*D603"

$ws.Range("A604").Value = 45986
$ws.Range("B604").Value = 'Amir Etien'
$ws.Range("C604").Value = 70
$ws.Range("D604").Value = 5
$ws.Range("E604").Value = 5
$ws.Range("F604").Value = 0
$ws.Range("H604").Value = 5
$ws.Range("I604").Formula = "=C604*D604"

$ws.Range("A605").Value = 45986
$ws.Range("B605").Value = 'Omar Benyounes'
$ws.Range("C605").Value = 70
$ws.Range("D605").Value = 5
$ws.Range("E605").Value = 6
$ws.Range("F605").Value = 3
$ws.Range("G605").Value = 'Ischio'
$ws.Range("H605").Value = 7
$ws.Range("I605").Formula = "=C605*D605"

$ws.Range("A606").Value = 45986
$ws.Range("B606").Value = 'Naim Ighbane'
$ws.Range("C606").Value = 70
$ws.Range("D606").Value = 3
$ws.Range("E606").Value = 6
$ws.Range("F606").Value = 5
$ws.Range("G606").Value = 'Genou'
$ws.Range("H606").Value = 4
$ws.Range("I606").Formula = "=C606*D606"

$ws.Range("A607").Value = 45986
$ws.Range("B607").Value = 'Karim Belmahi'
$ws.Range("C607").Value = 70
$ws.Range("D607").Value = 6
$ws.Range("E607").Value = 4
$ws.Range("F607").Value = 3
$ws.Range("G607").Value = 'Ischio'
$ws.Range("H607").Value = 10
$ws.Range("I607").Formula = "=C607*D607"

$ws.Range("A608").Value = 45986
$ws.Range("B608").Value = 'Maé Clavel'
$ws.Range("C608").Value = 70
$ws.Range("D608").Value = 5
$ws.Range("E608").Value = 1
$ws.Range("F608").Value = 0
$ws.Range("H608").Value = 7
$ws.Range("I608").Formula = "=C608*D608"

$ws.Range("A609").Value = 45986
$ws.Range("B609").Value = 'Jeremie Laurent'
$ws.Range("C609").Value = 70
$ws.Range("D609").Value = 7
$ws.Range("E609").Value = 3
$ws.Range("F609").Value = 0
$ws.Range("H609").Value = 9
$ws.Range("I609").Formula = "=C609*D609"

$ws.Range("A610").Value = 45986
$ws.Range("B610").Value = 'Levy Ndoutoume'
$ws.Range("C610").Value = 70
$ws.Range("D610").Value = 6
$ws.Range("E610").Value = 6
$ws.Range("F610").Value = 1
$ws.Range("G610").Value = 'Ischio'
$ws.Range("H610").Value = 8
$ws.Range("I610").Formula = "=C610*D610"

$ws.Range("A611").Value = 45986
$ws.Range("B611").Value = 'Hedi Nasri'
$ws.Range("C611").Value = 70
$ws.Range("D611").Value = 6
$ws.Range("E611").Value = 8
$ws.Range("F611").Value = 4
$ws.Range("G611").Value = 'Ischio courbature'
$ws.Range("H611").Value = 7
$ws.Range("I611").Formula = "=C611*D611"

$ws.Range("A612").Value = 45986
$ws.Range("B612").Value = 'Ilan Ihaddadene'
$ws.Range("C612").Value = 70
$ws.Range("D612").Value = 6
$ws.Range("E612").Value = 3
$ws.Range("F612").Value = 0
$ws.Range("H612").Value = 7
$ws.Range("I612").Formula = "=C612*D612"

$ws.Range("A613").Value = 45986
$ws.Range("B613").Value = 'Mattheo Haon'
$ws.Range("C613").Value = 70
$ws.Range("D613").Value = 5
$ws.Range("E613").Value = 0
$ws.Range("F613").Value = 0
$ws.Range("H613").Value = 10
$ws.Range("I613").Formula = "=C613*D613"

$ws.Range("A614").Value = 45986
$ws.Range("B614").Value = 'Karahali Souaré'
$ws.Range("C614").Value = 70
$ws.Range("D614").Value = 3
$ws.Range("E614").Value = 3
$ws.Range("F614").Value = 6
$ws.Range("G614").Value = 'Cheville '
$ws.Range("H614").Value = 3
$ws.Range("I614").Formula = "=C614*D614"

# Update the active view to match the saved selection
$ws.Range("L607").Select() | Out-Null

$excel.Calculate()
